# Remove needless imports on org.dozer
#
# Sheet1 has an "Environment" / "import" block in rows 7-9 listing the
# imports used by the mapping (org.openl.rules.mapping and
# org.openl.rules.mapping.to). Row 10 held an extra, no-longer-needed
# import entry ("org.dozer"). Deleting the whole row removes both the
# cell and its shared string, and shifts everything below it up by one
# row (merged ranges, the later data table, dimension, etc. all follow
# automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Rows("10:10").Delete()

# Make Sheet1 the active sheet/selection (it was Sheet2 before), with
# the cursor landing where the deleted row used to be.
$ws.Activate()
$ws.Range("D10").Select()
